$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 12592676
$ws.Range("I9").Value = 11111211
$ws.Range("K9").Value = 11111211
$ws.Range("M9").Value = -11111042

$ws.Range("H100").Value = 62501700
$ws.Range("I100").Value = 1967.5
$ws.Range("J100").Value = 250000900
$ws.Range("K100").Value = 1967.5
$ws.Range("L100").Value = 250000900
$ws.Range("M100").Value = -1426.5
$ws.Range("N100").Value = -250001982

$ws.Range("H113").Value = 3097.44
$ws.Range("I113").Value = 1853.125
$ws.Range("J113").Value = 5309.5557
$ws.Range("K113").Value = 1853.125
$ws.Range("L113").Value = 5309.5557
$ws.Range("M113").Value = 1400.875
$ws.Range("N113").Value = -11817.5557

$ws.Range("H137").Value = 1532.6727
$ws.Range("I137").Value = 1103.919
$ws.Range("J137").Value = 2414
$ws.Range("K137").Value = 3311.757000000001
$ws.Range("L137").Value = 7242
$ws.Range("M137").Value = -761.7570000000005
$ws.Range("N137").Value = -12342

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 34649
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H55").Value = 18150.857
$ws.Range("J55").Value = 18150.857
$ws.Range("L55").Value = 18150.857
$ws.Range("N55").Value = -18780.857

$ws.Range("H74").Value = 1023
$ws.Range("I74").Value = 1021.2
$ws.Range("K74").Value = 1021.2
$ws.Range("M74").Value = -147.2

$ws.Range("H77").Value = 1023
$ws.Range("I77").Value = 1021.2
$ws.Range("K77").Value = 5106
$ws.Range("M77").Value = -738

$ws.Range("H80").Value = 25050.8
$ws.Range("J80").Value = 25050.8
$ws.Range("L80").Value = 25050.8
$ws.Range("N80").Value = -27046.8

$ws.Range("H83").Value = 25050.8
$ws.Range("J83").Value = 25050.8
$ws.Range("L83").Value = 75152.39999999999
$ws.Range("N83").Value = -85136.39999999999

$ws.Range("H122").Value = 1294.8846
$ws.Range("I122").Value = 1291.2667
$ws.Range("J122").Value = 1299.8182
$ws.Range("K122").Value = 3873.800099999999
$ws.Range("L122").Value = 3899.4546
$ws.Range("M122").Value = -1423.800099999999
$ws.Range("N122").Value = -8799.454600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 33354.8
$ws.Range("J35").Value = 33354.8
$ws.Range("L35").Value = 33354.8
$ws.Range("N35").Value = -33974.8

$ws.Range("H82").Value = 23711.523
$ws.Range("I82").Value = 12094.25
$ws.Range("J82").Value = 26445
$ws.Range("K82").Value = 12094.25
$ws.Range("L82").Value = 26445
$ws.Range("M82").Value = -11711.25
$ws.Range("N82").Value = -27211

$ws.Range("H85").Value = 23711.523
$ws.Range("I85").Value = 12094.25
$ws.Range("J85").Value = 26445
$ws.Range("K85").Value = 12094.25
$ws.Range("L85").Value = 26445
$ws.Range("M85").Value = -10768.25
$ws.Range("N85").Value = -29097

$ws.Range("H99").Value = 2225.4546
$ws.Range("I99").Value = 1363.3334
$ws.Range("K99").Value = 1363.3334
$ws.Range("M99").Value = 134.6666

$ws.Range("H105").Value = 3195.3635
$ws.Range("I105").Value = 2127.6667
$ws.Range("J105").Value = 8000
$ws.Range("K105").Value = 2127.6667
$ws.Range("L105").Value = 8000
$ws.Range("M105").Value = -380.6667000000002
$ws.Range("N105").Value = -11494

$ws.Range("H130").Value = 39203.637
$ws.Range("J130").Value = 39203.637
$ws.Range("L130").Value = 39203.637
$ws.Range("N130").Value = -49243.637

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2479.352
$ws.Range("I31").Value = 1753.4615
$ws.Range("J31").Value = 3364.0312
$ws.Range("K31").Value = 1753.4615
$ws.Range("L31").Value = 3364.0312
$ws.Range("M31").Value = -1458.4615
$ws.Range("N31").Value = -3954.0312

$ws.Range("H34").Value = 2479.352
$ws.Range("I34").Value = 1753.4615
$ws.Range("J34").Value = 3364.0312
$ws.Range("K34").Value = 1753.4615
$ws.Range("L34").Value = 3364.0312
$ws.Range("M34").Value = -1551.4615
$ws.Range("N34").Value = -3768.0312

$ws.Range("H41").Value = 12850
$ws.Range("J41").Value = 19616.666
$ws.Range("L41").Value = 19616.666
$ws.Range("N41").Value = -20472.666

$ws.Range("H50").Value = 9101.833000000001
$ws.Range("J50").Value = 9101.833000000001
$ws.Range("L50").Value = 9101.833000000001
$ws.Range("N50").Value = -10351.833

$ws.Range("H59").Value = 16724.4
$ws.Range("J59").Value = 16724.4
$ws.Range("L59").Value = 16724.4
$ws.Range("N59").Value = -19014.4

$ws.Range("H105").Value = 4041.818
$ws.Range("I105").Value = 1200
$ws.Range("J105").Value = 6410
$ws.Range("K105").Value = 1200
$ws.Range("L105").Value = 6410
$ws.Range("M105").Value = 547
$ws.Range("N105").Value = -9904

$ws.Range("H109").Value = 10828.571
$ws.Range("J109").Value = 10828.571
$ws.Range("L109").Value = 10828.571
$ws.Range("N109").Value = -12908.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4250.75
$ws.Range("I68").Value = 7751
$ws.Range("J68").Value = 750.5
$ws.Range("K68").Value = 23253
$ws.Range("L68").Value = 2251.5
$ws.Range("M68").Value = -22442
$ws.Range("N68").Value = -3873.5

$ws.Range("H71").Value = 4250.75
$ws.Range("I71").Value = 7751
$ws.Range("J71").Value = 750.5
$ws.Range("K71").Value = 69759
$ws.Range("L71").Value = 6754.5
$ws.Range("M71").Value = -65703
$ws.Range("N71").Value = -14866.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3017.8333
$ws.Range("I97").Value = 2808.6428
$ws.Range("J97").Value = 3750
$ws.Range("K97").Value = 2808.6428
$ws.Range("L97").Value = 3750
$ws.Range("M97").Value = -2312.6428
$ws.Range("N97").Value = -4742

$ws.Range("H122").Value = 1891.5454
$ws.Range("I122").Value = 1788.375
$ws.Range("K122").Value = 5365.125
$ws.Range("M122").Value = -2915.125

$ws.Range("H123").Value = 23103.715
$ws.Range("J123").Value = 23103.715
$ws.Range("L123").Value = 23103.715
$ws.Range("N123").Value = -28003.715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1198.2424
$ws.Range("I68").Value = 1055.3684
$ws.Range("J68").Value = 1392.1428
$ws.Range("K68").Value = 1055.3684
$ws.Range("L68").Value = 1392.1428
$ws.Range("M68").Value = -306.3684000000001
$ws.Range("N68").Value = -2890.1428

$ws.Range("H71").Value = 1198.2424
$ws.Range("I71").Value = 1055.3684
$ws.Range("J71").Value = 1392.1428
$ws.Range("K71").Value = 5276.842000000001
$ws.Range("L71").Value = 6960.714
$ws.Range("M71").Value = -1532.842000000001
$ws.Range("N71").Value = -14448.714

$ws.Range("H100").Value = 1475.3334
$ws.Range("I100").Value = 1545
$ws.Range("J100").Value = 1336
$ws.Range("K100").Value = 1545
$ws.Range("L100").Value = 1336
$ws.Range("M100").Value = -1004
$ws.Range("N100").Value = -2418

$ws.Range("H122").Value = 4128.625
$ws.Range("I122").Value = 2425
$ws.Range("J122").Value = 4696.5
$ws.Range("K122").Value = 7275
$ws.Range("L122").Value = 14089.5
$ws.Range("M122").Value = -4825
$ws.Range("N122").Value = -18989.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws.Range("H107").Value = 1630.9459
$ws.Range("I107").Value = 971.2857
$ws.Range("J107").Value = 3683.2222
$ws.Range("K107").Value = 2913.8571
$ws.Range("L107").Value = 11049.6666
$ws.Range("M107").Value = -993.8571000000002
$ws.Range("N107").Value = -14889.6666

$ws.Range("H109").Value = 17022.334
$ws.Range("J109").Value = 17022.334
$ws.Range("L109").Value = 17022.334
$ws.Range("N109").Value = -19796.334
